# Update TPM-derived statistics in the LR-pairs sheet (Cp-Slc40a1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Range("G2").Value = 3.702073333333333
$ws.Range("H2").Value = 11.10622
$ws.Range("I2").Value = 0.031699224716142
$ws.Range("J2").Value = 0.031699224716142
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.05651999999999999
$ws.Range("N2").Value = 0.16956
$ws.Range("O2").Value = 0.01567922357698054
$ws.Range("P2").Value = 0.01567922357698054
$ws.Range("Q2").Value = 0.2092411848
$ws.Range("R2").Value = 1.8831706632
$ws.Range("S2").Value = 0.000497019231541338
$ws.Range("T2").Value = 0.0004970192315413379

# Row 3 (ECs -> FAPs)
$ws.Range("G3").Value = 3.702073333333333
$ws.Range("H3").Value = 11.10622
$ws.Range("I3").Value = 0.031699224716142
$ws.Range("J3").Value = 0.031699224716142
$ws.Range("O3").Value = 0.4433226490342289
$ws.Range("P3").Value = 0.4433226490342288
$ws.Range("Q3").Value = 5.916195778264443
$ws.Range("R3").Value = 53.24576200437999
$ws.Range("S3").Value = 0.01405298427349137
$ws.Range("T3").Value = 0.01405298427349137

# Row 4 (ECs -> MuSCs)
$ws.Range("G4").Value = 3.702073333333333
$ws.Range("H4").Value = 11.10622
$ws.Range("I4").Value = 0.031699224716142
$ws.Range("J4").Value = 0.031699224716142
$ws.Range("M4").Value = 1.950174
$ws.Range("N4").Value = 5.850522
$ws.Range("O4").Value = 0.5409981273887907
$ws.Range("P4").Value = 0.5409981273887906
$ws.Range("Q4").Value = 7.219687160759999
$ws.Range("R4").Value = 64.97718444683998
$ws.Range("S4").Value = 0.01714922121110929
$ws.Range("T4").Value = 0.01714922121110929

# Row 5 (FAPs -> ECs)
$ws.Range("I5").Value = 0.3975581268808804
$ws.Range("J5").Value = 0.3975581268808804
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.05651999999999999
$ws.Range("N5").Value = 0.16956
$ws.Range("O5").Value = 0.01567922357698054
$ws.Range("P5").Value = 0.01567922357698054
$ws.Range("Q5").Value = 2.62421350176
$ws.Range("R5").Value = 23.61792151584
$ws.Range("S5").Value = 0.006233402756210923
$ws.Range("T5").Value = 0.006233402756210921

# Row 6 (FAPs -> FAPs)
$ws.Range("I6").Value = 0.3975581268808804
$ws.Range("J6").Value = 0.3975581268808804
$ws.Range("O6").Value = 0.4433226490342289
$ws.Range("P6").Value = 0.4433226490342288
$ws.Range("S6").Value = 0.176246521953918
$ws.Range("T6").Value = 0.176246521953918

# Row 7 (FAPs -> MuSCs)
$ws.Range("I7").Value = 0.3975581268808804
$ws.Range("J7").Value = 0.3975581268808804
$ws.Range("M7").Value = 1.950174
$ws.Range("N7").Value = 5.850522
$ws.Range("O7").Value = 0.5409981273887907
$ws.Range("P7").Value = 0.5409981273887906
$ws.Range("Q7").Value = 90.54623038891199
$ws.Range("R7").Value = 814.916073500208
$ws.Range("S7").Value = 0.2150782021707516
$ws.Range("T7").Value = 0.2150782021707515

# Row 8 (MuSCs -> ECs)
$ws.Range("G8").Value = 66.65560933333335
$ws.Range("H8").Value = 199.966828
$ws.Range("I8").Value = 0.5707426484029776
$ws.Range("J8").Value = 0.5707426484029775
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.05651999999999999
$ws.Range("N8").Value = 0.16956
$ws.Range("O8").Value = 0.01567922357698054
$ws.Range("P8").Value = 0.01567922357698054
$ws.Range("Q8").Value = 3.76737503952
$ws.Range("R8").Value = 33.90637535568
$ws.Range("S8").Value = 0.008948801589228284
$ws.Range("T8").Value = 0.008948801589228278

# Row 9 (MuSCs -> FAPs)
$ws.Range("G9").Value = 66.65560933333335
$ws.Range("H9").Value = 199.966828
$ws.Range("I9").Value = 0.5707426484029776
$ws.Range("J9").Value = 0.5707426484029775
$ws.Range("O9").Value = 0.4433226490342289
$ws.Range("P9").Value = 0.4433226490342288
$ws.Range("Q9").Value = 106.5207517595125
$ws.Range("R9").Value = 958.6867658356121
$ws.Range("S9").Value = 0.2530231428068195
$ws.Range("T9").Value = 0.2530231428068194

# Row 10 (MuSCs -> MuSCs)
$ws.Range("G10").Value = 66.65560933333335
$ws.Range("H10").Value = 199.966828
$ws.Range("I10").Value = 0.5707426484029776
$ws.Range("J10").Value = 0.5707426484029775
$ws.Range("M10").Value = 1.950174
$ws.Range("N10").Value = 5.850522
$ws.Range("O10").Value = 0.5409981273887907
$ws.Range("P10").Value = 0.5409981273887906
$ws.Range("Q10").Value = 129.990036276024
$ws.Range("R10").Value = 1169.910326484216
$ws.Range("S10").Value = 0.3087707040069299
$ws.Range("T10").Value = 0.3087707040069297
